# Generate Report for Handoff
# This script updates the localization-status workbook to reflect that the
# 968b3c4b... file is now "Ready for handoff" (with new handoff timestamps),
# and removes the bc4e5720... file row, which is no longer part of the report.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

# Row 2 (968b3c4b...): Status changes from "Handed back: in sync with en-US"
# to "Ready for handoff"
$ws1.Range("B2").Value = "Ready for handoff"
$ws1.Range("C2").Value = "Ready for handoff"

# Row 3 (bc4e5720...) is removed entirely; remaining rows shift up
$ws1.Rows.Item(3).Delete()

# Fix up the hyperlinks collection to reflect the new layout
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/34442f09a1faea49bb8db58846ca872e779e2a59/e2e/968b3c4b-617e-460b-ac64-45d96ecf4a67.md", "", "", "968b3c4b-617e-460b-ac64-45d96ecf4a67.md")
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/34442f09a1faea49bb8db58846ca872e779e2a59/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

# Row 2 (968b3c4b...): Status + Latest Handoff Datetime change
$ws2.Range("B2").Value = "Ready for handoff"
$ws2.Range("D2").Value = "2016-03-02 15:06:11"

# Row 3 (bc4e5720...) is removed entirely; remaining rows shift up
$ws2.Rows.Item(3).Delete()

# Fix up the hyperlinks collection to reflect the new layout
$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/34442f09a1faea49bb8db58846ca872e779e2a59/e2e/968b3c4b-617e-460b-ac64-45d96ecf4a67.md", "", "", "968b3c4b-617e-460b-ac64-45d96ecf4a67.md")
$ws2.Hyperlinks.Add($ws2.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/67aa3f79304f991214dd8abd5965ffc74580436b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/968b3c4b-617e-460b-ac64-45d96ecf4a67.8341d5b57efefc4be07c960d513f43ef16b5973c.zh-cn.xlf", "", "", "968b3c4b-617e-460b-ac64-45d96ecf4a67.8341d5b57efefc4be07c960d513f43ef16b5973c.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/bf12176fe7cdac8bfcd46c435ed8232259580438/e2e/968b3c4b-617e-460b-ac64-45d96ecf4a67.md", "", "", "968b3c4b-617e-460b-ac64-45d96ecf4a67.md")
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/873fa4cb86a6157b14688f274c0f932e8a8c2b74/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/968b3c4b-617e-460b-ac64-45d96ecf4a67.8341d5b57efefc4be07c960d513f43ef16b5973c.zh-cn.xlf", "", "", "968b3c4b-617e-460b-ac64-45d96ecf4a67.8341d5b57efefc4be07c960d513f43ef16b5973c.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/34442f09a1faea49bb8db58846ca872e779e2a59/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

# Row 2 (968b3c4b...): Status + Latest Handoff Datetime change
$ws3.Range("B2").Value = "Ready for handoff"
$ws3.Range("D2").Value = "2016-03-02 15:06:26"

# Row 3 (bc4e5720...) is removed entirely; remaining rows shift up
$ws3.Rows.Item(3).Delete()

# Fix up the hyperlinks collection to reflect the new layout
$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/34442f09a1faea49bb8db58846ca872e779e2a59/e2e/968b3c4b-617e-460b-ac64-45d96ecf4a67.md", "", "", "968b3c4b-617e-460b-ac64-45d96ecf4a67.md")
$ws3.Hyperlinks.Add($ws3.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dbf3fdcac76906e9d9c8d0cc33672277287e0683/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/968b3c4b-617e-460b-ac64-45d96ecf4a67.8341d5b57efefc4be07c960d513f43ef16b5973c.de-de.xlf", "", "", "968b3c4b-617e-460b-ac64-45d96ecf4a67.8341d5b57efefc4be07c960d513f43ef16b5973c.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/ad457bcc16e319d37b33cf6ace6ae1f3cd779b40/e2e/968b3c4b-617e-460b-ac64-45d96ecf4a67.md", "", "", "968b3c4b-617e-460b-ac64-45d96ecf4a67.md")
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/c8529ce38153b6b7a752d2889a03275c9ea6ab1e/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/968b3c4b-617e-460b-ac64-45d96ecf4a67.8341d5b57efefc4be07c960d513f43ef16b5973c.de-de.xlf", "", "", "968b3c4b-617e-460b-ac64-45d96ecf4a67.8341d5b57efefc4be07c960d513f43ef16b5973c.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/34442f09a1faea49bb8db58846ca872e779e2a59/.localization-config", "", "", ".localization-config")
